# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.274.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").Value = "'3.321.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("D5").Value = "'592.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").Value = "'187.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.45%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'0.608"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.04%  "

$ws.Range("D9").Value = "'0.137"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.55%  "

$ws.Range("D10").Value = "'6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").Value = "'0.426"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("D12").Value = "'3.891.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.27%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").Value = "'29.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.94%  "

$ws.Range("D15").Value = "'69.251.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.62%  "

$ws.Range("D16").Value = "'0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.23%  "

$ws.Range("D17").Value = "'3.298.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").Value = "'5.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").Value = "'13.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.59%  "

$ws.Range("D20").Value = "'390.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.46%  "

$ws.Range("D21").Value = "'7.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.79%  "

$ws.Range("D22").Value = "'71.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  +4.32%  "

$ws.Range("D25").Value = "'0.524"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.16%  "

$ws.Range("E26").Value = "  +6.12%  "

$ws.Range("D27").Value = "'9.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.16%  "

$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").Value = "'5.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.91%  "

$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("E31").Value = "  +5.44%  "

$ws.Range("D32").Value = "'23.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.12%  "

$ws.Range("D33").Value = "'7.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.70%  "

$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "'1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.14%  "

$ws.Range("D36").Value = "'163.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.68%  "

$ws.Range("D37").Value = "'1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.59%  "

$ws.Range("D38").Value = "'0.845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "

$ws.Range("D39").Value = "'27.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.29%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.33%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.91%  "

$ws.Range("D43").Value = "'26.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.94%  "

$ws.Range("D44").Value = "'0.0702"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.38%  "

$ws.Range("D45").Value = "'41.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.74%  "

$ws.Range("D46").Value = "'2.670.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("D47").Value = "'344.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.78%  "

$ws.Range("D48").Value = "'0.0288"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.73%  "

$ws.Range("D49").Value = "'32.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.10%  "

$ws.Range("D50").Value = "'1.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.31%  "
